# ENVS363_563_Redesign.xlsx -- "more structure for website"
#
# 1. Sessions_TODOs: add a course-website link, an "Other To Dos" section
#    and a "Demonstratos" (demonstrators) roster below the existing table.
# 2. Learning_Outcomes: clear the (invisible/no-op) border formatting that
#    was sitting on the whole "Assessment weighting" column.
# 3. Add a new "random" sheet at the end of the workbook with a scratch
#    link on it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Sessions_TODOs
# ---------------------------------------------------------------------
$todo = $wb.Worksheets.Item("Sessions_TODOs")
$todo.Activate()

$todo.Range("A16").Value = "https://darribas.org/gds_course/content/home.html"

$todo.Range("A18").Value = "Other To Dos"

$todo.Range("A19").Value = "get a better website name"
$todo.Range("B19").Value = "Eli"
$todo.Range("B11").Copy()
$todo.Range("B19").PasteSpecial(-4122)

$todo.Range("A21").Value = "Demonstratos"

$todo.Range("A22").Value = "Matthew Howard"
$todo.Range("A23").Value = "Ruth Neville"
$todo.Range("A24").Value = "Akos Balog"
$todo.Range("A25").Value = "Rodgers Iradukunda"
$todo.Range("D9").Copy()
$todo.Range("A22:A25").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# widen column A so the new link/notes are readable
$todo.Columns.Item(1).ColumnWidth = 26.67

$todo.Range("A16").Select()

# ---------------------------------------------------------------------
# 2. Learning_Outcomes -- drop the stray border formatting on column C
# ---------------------------------------------------------------------
$lo = $wb.Worksheets.Item("Learning_Outcomes")
$lo.Range("C1:C6").Borders.LineStyle = -4142

# ---------------------------------------------------------------------
# 3. New "random" scratch sheet at the end
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$randomSheet = $wb.Worksheets.Add($null, $lastSheet)
$randomSheet.Name = "random"
$randomSheet.Range("A1").Value = "https://github.com/hadley/r4ds/blob/main/_quarto.yml"

# restore focus to the main tracking sheet
$todo.Activate()
$todo.Range("A16").Select()
